# Add a new "LIMITED COMPANY" client type block to both the
# "Applicability" and "Fees" sheets, mirroring the existing 34-row
# Service/SubService pattern used for every other ClientType.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Applicability")
$ws2 = $wb.Worksheets.Item("Fees")

$clientType = "LIMITED COMPANY"

$rows = @(
    @{ Row=342; Service='Incorporation'; SubService=$null; Applicable=$true; Fee=25000 },
    @{ Row=343; Service='Accounting'; SubService='Monthly Accounting'; Applicable=$true; Fee=150000 },
    @{ Row=344; Service='Accounting'; SubService='Quarterly Accounting'; Applicable=$true; Fee=100000 },
    @{ Row=345; Service='Accounting'; SubService='Half Yearly Accounting'; Applicable=$true; Fee=50000 },
    @{ Row=346; Service='Accounting'; SubService='Annual Accounting'; Applicable=$true; Fee=25000 },
    @{ Row=347; Service='Filing of Income Tax Returns'; SubService=$null; Applicable=$true; Fee=7500 },
    @{ Row=348; Service='GST Registration'; SubService=$null; Applicable=$true; Fee=5000 },
    @{ Row=349; Service='Filing of GSTR Returns'; SubService='Filing of GSTR 1 Return'; Applicable=$true; Fee=12000 },
    @{ Row=350; Service='Filing of GSTR Returns'; SubService='Filing of GSTR 3B Return'; Applicable=$true; Fee=12000 },
    @{ Row=351; Service='Filing of GSTR Returns'; SubService='Filing of GSTR 9'; Applicable=$true; Fee=30000 },
    @{ Row=352; Service='Filing of GSTR Returns'; SubService='Filing of GSTR 9C'; Applicable=$true; Fee=30000 },
    @{ Row=353; Service='Statutory Audit'; SubService=$null; Applicable=$true; Fee=25000 },
    @{ Row=354; Service='Income Tax Audit'; SubService=$null; Applicable=$true; Fee=25000 },
    @{ Row=355; Service='ROC Filing'; SubService='Filing of Form 8'; Applicable=$false; Fee=0 },
    @{ Row=356; Service='ROC Filing'; SubService='Filing of Form 11'; Applicable=$false; Fee=0 },
    @{ Row=357; Service='ROC Filing'; SubService='Filing of Form AOC 4'; Applicable=$true; Fee=5000 },
    @{ Row=358; Service='ROC Filing'; SubService='Filing of Form MGT 7'; Applicable=$true; Fee=5000 },
    @{ Row=359; Service='TDS Return'; SubService='Filing of TDS Return in Form 26Q'; Applicable=$true; Fee=10000 },
    @{ Row=360; Service='TDS Return'; SubService='Filing of TDS Return in Form 24Q'; Applicable=$true; Fee=10000 },
    @{ Row=361; Service='TDS Return'; SubService='Filing of TDS Return in Form 27Q'; Applicable=$true; Fee=10000 },
    @{ Row=362; Service='TDS Return'; SubService='Filing of TDS Return in Form 26QB'; Applicable=$true; Fee=2500 },
    @{ Row=363; Service='TDS Return'; SubService='Filing of TDS Return in Form 26QC'; Applicable=$true; Fee=2500 },
    @{ Row=364; Service='Profession Tax Registration'; SubService='PTEC Registration'; Applicable=$true; Fee=2000 },
    @{ Row=365; Service='Profession Tax Registration'; SubService='PTRC Registration'; Applicable=$true; Fee=2000 },
    @{ Row=366; Service='Profession Tax Returns'; SubService='Monthly PTRC Return Filing'; Applicable=$true; Fee=6000 },
    @{ Row=367; Service='Profession Tax Returns'; SubService='Annual PTRC Return Filing'; Applicable=$true; Fee=2500 },
    @{ Row=368; Service='Annual PTEC Payment'; SubService=$null; Applicable=$true; Fee=500 },
    @{ Row=369; Service='Event Based Filing'; SubService='DIR 12'; Applicable=$true; Fee=2500 },
    @{ Row=370; Service='Event Based Filing'; SubService='ADT 1'; Applicable=$true; Fee=2500 },
    @{ Row=371; Service='Event Based Filing'; SubService='ADT 3'; Applicable=$true; Fee=2500 },
    @{ Row=372; Service='Event Based Filing'; SubService='Change of Address in ROC'; Applicable=$true; Fee=2500 },
    @{ Row=373; Service='Event Based Filing'; SubService='Change of Address in GST'; Applicable=$true; Fee=2500 },
    @{ Row=374; Service='Event Based Filing'; SubService='ROC E-Kyc for Directors'; Applicable=$true; Fee=1500 },
    @{ Row=375; Service='Event Based Filing'; SubService='MSME Application'; Applicable=$true; Fee=2500 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # --- Applicability sheet: Service / SubService / ClientType / Applicable (bool) ---
    $ws1.Cells.Item($row, 1).Value = $r.Service
    if ($r.SubService -ne $null) {
        $ws1.Cells.Item($row, 2).Value = $r.SubService
    }
    $ws1.Cells.Item($row, 3).Value = $clientType
    $ws1.Cells.Item($row, 4).Value = $r.Applicable

    # --- Fees sheet: Service / SubService / ClientType / FeeINR (number) ---
    $ws2.Cells.Item($row, 1).Value = $r.Service
    if ($r.SubService -ne $null) {
        $ws2.Cells.Item($row, 2).Value = $r.SubService
    }
    $ws2.Cells.Item($row, 3).Value = $clientType
    $ws2.Cells.Item($row, 4).Value = $r.Fee
}

# Fees sheet gains an explicit column-width block (bestFit-style) matching
# the Applicability sheet's A:C widths, with its own (narrower) D width.
$ws2.Columns.Item(1).ColumnWidth = 25.59
$ws2.Columns.Item(2).ColumnWidth = 30.31
$ws2.Columns.Item(3).ColumnWidth = 17.45
$ws2.Columns.Item(4).ColumnWidth = 6.59

# Move the Fees sheet's view down to the newly-added block and select C342,
# matching where the edit was made.
$ws2.Activate()
$ws2.Range("C342").Select()
$win = $excel.Windows.Item(1)
$win.ScrollRow = 332
$win.ScrollColumn = 1

# Restore Applicability as the active/visible tab.
$ws1.Activate()
